$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 144 - date 2025-11-11, station index 26 (四方坪站充电量(kw))
$ws.Cells.Item(144, 1).Value = "2025-11-11"
$ws.Cells.Item(144, 2).Value = "四方坪站充电量(kw)"
$row144 = @(492.38900000000001, 1235.8040000000005, 449.37, 405.92999999999995, 354.09699999999998, 520.25099999999998, 561.40199999999993, 214.316, 54.720000000000006, 47.040000000000006, 243.93400000000003, 157.69, 613.20900000000006, 1440.0790000000002, 791.76699999999983, 173.71000000000004, 381.56200000000001, 183.21700000000004, 36.82, 83.86, 90.64, 35.81, 87.94, 60.860999999999997)
for ($i = 0; $i -lt $row144.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(144, $col).Value = $row144[$i]
}

# New row 145 - date 2025-11-11, station index 27 (高岭站充电量(kw))
$ws.Cells.Item(145, 1).Value = "2025-11-11"
$ws.Cells.Item(145, 2).Value = "高岭站充电量(kw)"
$row145 = @(100.717, 331.10899999999998, 93.274000000000001, 32.695999999999998, 34.195, 187.80799999999999, 212.44699999999997, 244.47299999999998, 300.30399999999997, 170.255, 135.40199999999999, 158.13899999999998, 421.839, 656.18600000000004, 295.87199999999996, 230.15099999999998, 175.36799999999999, 135.42599999999999, 46.097000000000001, 38.814999999999998, 107.387, 0, 0, 17.109000000000002)
for ($i = 0; $i -lt $row145.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(145, $col).Value = $row145[$i]
}

# Apply the number formats matching the rest of the column (date + 0.00 style)
$ws.Range("A144:A145").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C144:Z145").NumberFormat = "0.00_);[Red]\(0.00\)"

# Update selection to match the final cursor position recorded in the edit
$null = $ws.Range("E148").Select()
